$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 182 previously had blank "Hora de Reparacion" / "Tiempo de Reparacion"
# placeholder cells (F182/G182). Clear them so they no longer appear at all,
# matching the other "no repair time recorded" rows elsewhere in the sheet.
$ws.Range("F182").ClearContents()
$ws.Range("G182").ClearContents()

# New incident row 183
$ws.Range("A183").Value = "SPL"
$ws.Range("B183").Value = "No detecta marcas Power"
$ws.Range("C183").Value = "'2024-06-10"
$ws.Range("D183").Value = "12:04:52"
$ws.Range("E183").Value = "Mañana"
$ws.Range("F183").Value = "12:04:53"
$ws.Range("G183").Value = "0:00:01"
$ws.Range("H183").Value = "-0.00 minutos"

# New incident row 184
$ws.Range("A184").Value = "WC47 NACP"
$ws.Range("B184").Value = "Fallo en elevador"
$ws.Range("C184").Value = "'2024-06-10"
$ws.Range("D184").Value = "12:08:34"
$ws.Range("E184").Value = "Mañana"
$ws.Range("F184").Value = "12:08:36"
$ws.Range("G184").Value = "0:00:02"
$ws.Range("H184").Value = "-0.01 minutos"

# New incident row 185
$ws.Range("A185").Value = "WC47 NACP"
$ws.Range("B185").Value = "Fallo en paletizador"
$ws.Range("C185").Value = "'2024-06-10"
$ws.Range("D185").Value = "12:13:59"
$ws.Range("E185").Value = "Mañana"
$ws.Range("F185").Value = "12:14:00"
$ws.Range("G185").Value = "0:00:01"
$ws.Range("H185").Value = "-0.01 minutos"
